# Apply crypto price/symbol-list refresh per the commit diff.
# Columns B/C/E hold non-numeric text -> assign directly.
# Columns D/G hold numeric-looking text (prices/hour) that must stay
# TEXT (not auto-converted to Number) so formatting such as trailing
# zeros is preserved exactly as in the source diff -> pre-format the
# cell as Text ("@") before writing the value, exactly as a user would
# need to do in real Excel to keep a numeric-looking entry as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "254.07"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "10"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.93"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "10"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.118"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "10"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06026"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "10"

# Row 6
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.608"
$ws.Range("E6").Value = "5KuCoinTokenKCS"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "10"

# Row 7
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.437"
$ws.Range("E7").Value = "6GateTokenGT"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "10"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.323"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "10"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8013"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "10"

# Row 10
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01279"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "10"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1524"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "10"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07959"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "10"

# Row 13
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03352"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "10"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03097"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "10"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09305"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "10"

# Row 16
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.609"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "10"

# Row 17
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001676"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "10"

# Row 18
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04790"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "10"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006267"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "10"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005898"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "10"

# Row 21
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "10"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001507"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "10"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.690"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "10"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.200"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "10"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3347"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "10"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1268"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "10"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0006504"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "10"

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "10"

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "10"

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "10"

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "10"

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "10"

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "10"

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "10"

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "10"

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "10"

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "10"

# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "10"

# Row 39
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "10"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04484"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "10"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007046"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "10"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1074"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "10"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003365"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "10"

# Row 44
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "10"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.002472"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "10"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005898"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "10"

# Row 47
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "10"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7033"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "10"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.09460"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "10"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002110"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "10"

# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "10"
